$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 474.54544
$ws.Range("I19").Value = 345.42856
$ws.Range("J19").Value = 534.8
$ws.Range("K19").Value = 345.42856
$ws.Range("L19").Value = 534.8
$ws.Range("M19").Value = -170.42856
$ws.Range("N19").Value = -884.8

$ws.Range("H123").Value = 88000
$ws.Range("J123").Value = 88000
$ws.Range("L123").Value = 88000
$ws.Range("N123").Value = -97800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2419.7727
$ws.Range("I2").Value = 2438.3635
$ws.Range("J2").Value = 2401.182
$ws.Range("K2").Value = 2438.3635
$ws.Range("L2").Value = 2401.182
$ws.Range("M2").Value = -2325.3635
$ws.Range("N2").Value = -2627.182

$ws.Range("H32").Value = 5818.23
$ws.Range("I32").Value = 4659.1445
$ws.Range("J32").Value = 16250
$ws.Range("K32").Value = 4659.1445
$ws.Range("L32").Value = 16250
$ws.Range("M32").Value = -4372.1445
$ws.Range("N32").Value = -16824

$ws.Range("H82").Value = 33454
$ws.Range("J82").Value = 33454
$ws.Range("L82").Value = 33454
$ws.Range("N82").Value = -34176

$ws.Range("H85").Value = 33454
$ws.Range("J85").Value = 33454
$ws.Range("L85").Value = 33454
$ws.Range("N85").Value = -35950

$ws.Range("H116").Value = 2419.7727
$ws.Range("I116").Value = 2438.3635
$ws.Range("J116").Value = 2401.182
$ws.Range("K116").Value = 2438.3635
$ws.Range("L116").Value = 2401.182
$ws.Range("M116").Value = -144.3634999999999
$ws.Range("N116").Value = -6989.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2419.7727
$ws.Range("I3").Value = 2438.3635
$ws.Range("J3").Value = 2401.182
$ws.Range("K3").Value = 2438.3635
$ws.Range("L3").Value = 2401.182
$ws.Range("M3").Value = -2324.3635
$ws.Range("N3").Value = -2629.182

$ws.Range("H20").Value = 1678.7858
$ws.Range("I20").Value = 1468
$ws.Range("J20").Value = 1889.5714
$ws.Range("K20").Value = 1468
$ws.Range("L20").Value = 1889.5714
$ws.Range("M20").Value = -1221
$ws.Range("N20").Value = -2383.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 697.64105
$ws.Range("I5").Value = 336.12
$ws.Range("J5").Value = 1343.2142
$ws.Range("K5").Value = 1008.36
$ws.Range("L5").Value = 4029.6426
$ws.Range("M5").Value = -896.36
$ws.Range("N5").Value = -4253.642599999999

$ws.Range("H6").Value = 371.2857
$ws.Range("I6").Value = 59.8
$ws.Range("J6").Value = 1150
$ws.Range("K6").Value = 179.4
$ws.Range("L6").Value = 3450
$ws.Range("M6").Value = -66.39999999999998
$ws.Range("N6").Value = -3676

$ws.Range("H64").Value = 3581.75
$ws.Range("I64").Value = 1492.5714
$ws.Range("J64").Value = 4442
$ws.Range("K64").Value = 4477.7142
$ws.Range("L64").Value = 13326
$ws.Range("M64").Value = -4207.7142
$ws.Range("N64").Value = -13866

$ws.Range("H67").Value = 3581.75
$ws.Range("I67").Value = 1492.5714
$ws.Range("J67").Value = 4442
$ws.Range("K67").Value = 4477.7142
$ws.Range("L67").Value = 13326
$ws.Range("M67").Value = -3541.7142
$ws.Range("N67").Value = -15198

$ws.Range("H75").Value = 3266.6924
$ws.Range("I75").Value = 1710.8
$ws.Range("J75").Value = 4239.125
$ws.Range("K75").Value = 5132.4
$ws.Range("L75").Value = 12717.375
$ws.Range("M75").Value = -4134.4
$ws.Range("N75").Value = -14713.375

$ws.Range("H78").Value = 3266.6924
$ws.Range("I78").Value = 1710.8
$ws.Range("J78").Value = 4239.125
$ws.Range("K78").Value = 15397.2
$ws.Range("L78").Value = 38152.125
$ws.Range("M78").Value = -10405.2
$ws.Range("N78").Value = -48136.125

$ws.Range("H108").Value = 763.2857
$ws.Range("J108").Value = 3030
$ws.Range("L108").Value = 9090
$ws.Range("N108").Value = -14850

$ws.Range("H113").Value = 587.4706
$ws.Range("I113").Value = 552
$ws.Range("J113").Value = 622.94116
$ws.Range("K113").Value = 1656
$ws.Range("L113").Value = 1868.82348
$ws.Range("M113").Value = 514
$ws.Range("N113").Value = -6208.82348

$ws.Range("H131").Value = 897.6061
$ws.Range("I131").Value = 381.7143
$ws.Range("J131").Value = 1036.5
$ws.Range("K131").Value = 1145.1429
$ws.Range("L131").Value = 3109.5
$ws.Range("M131").Value = 3894.8571
$ws.Range("N131").Value = -13189.5

$ws.Range("H135").Value = 697.64105
$ws.Range("I135").Value = 336.12
$ws.Range("J135").Value = 1343.2142
$ws.Range("K135").Value = 3025.08
$ws.Range("L135").Value = 12088.9278
$ws.Range("M135").Value = -490.0799999999999
$ws.Range("N135").Value = -17158.9278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2175
$ws.Range("I113").Value = 1612
$ws.Range("J113").Value = 2738
$ws.Range("K113").Value = 1612
$ws.Range("L113").Value = 2738
$ws.Range("M113").Value = 558
$ws.Range("N113").Value = -7078

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 639.6087
$ws.Range("I22").Value = 509.9091
$ws.Range("J22").Value = 758.5
$ws.Range("K22").Value = 509.9091
$ws.Range("L22").Value = 758.5
$ws.Range("M22").Value = -214.9091
$ws.Range("N22").Value = -1348.5

$ws.Range("H27").Value = 639.6087
$ws.Range("I27").Value = 509.9091
$ws.Range("J27").Value = 758.5
$ws.Range("K27").Value = 509.9091
$ws.Range("L27").Value = 758.5
$ws.Range("M27").Value = -402.9091
$ws.Range("N27").Value = -972.5

$ws.Range("H68").Value = 1867
$ws.Range("I68").Value = 1800.5
$ws.Range("K68").Value = 1800.5
$ws.Range("M68").Value = -1051.5

$ws.Range("H71").Value = 1867
$ws.Range("I71").Value = 1800.5
$ws.Range("K71").Value = 9002.5
$ws.Range("M71").Value = -5258.5

$ws.Range("H132").Value = 31122.084
$ws.Range("I132").Value = 3463.2273
$ws.Range("J132").Value = 74586
$ws.Range("K132").Value = 10389.6819
$ws.Range("L132").Value = 223758
$ws.Range("M132").Value = -7859.6819
$ws.Range("N132").Value = -228818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

$ws.Range("H100").Value = 92163.55
$ws.Range("I100").Value = 83983.336
$ws.Range("J100").Value = 101979.8
$ws.Range("K100").Value = 167966.672
$ws.Range("L100").Value = 203959.6
$ws.Range("M100").Value = -167425.672
$ws.Range("N100").Value = -205041.6

$ws.Range("H140").Value = 53805
$ws.Range("J140").Value = 53805
$ws.Range("L140").Value = 53805
$ws.Range("N140").Value = -64165
